# Atualização de bases das ligas, do dia: 11-06-2024 às 21:19
#
# The underlying data rows for two match pairs were swapped in the source
# feed (the rows kept their sequential "id" in column A, but all the other
# match data - id (column B) through PL_AhUnder (column AD) - moved to the
# other row of the pair). Re-create that by swapping the B:AD contents of
# row 95 <-> row 96, and row 173 <-> row 174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($rowA, $rowB) {
    $rangeA = "B$rowA`:AD$rowA"
    $rangeB = "B$rowB`:AD$rowB"

    $valsA = $ws.Range($rangeA).Value2
    $valsB = $ws.Range($rangeB).Value2

    $ws.Range($rangeA).Value2 = $valsB
    $ws.Range($rangeB).Value2 = $valsA
}

Swap-RowData 95 96
Swap-RowData 173 174
